$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.049.62"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = "  -0.51%  "
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.828.25"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = "  -0.38%  "
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9992"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = "  -0.01%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.30"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +0.23%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6357"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -4.54%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  +0.02%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.82"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  +6.76%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2936"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  +0.58%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07332"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  -0.43%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22.81"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  +0.80%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07652"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  -0.56%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.828.93"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  -0.02%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.983"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  +0.21%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6637"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  -0.24%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "81.88"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  -1.91%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.052"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  -0.53%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008643"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  +4.63%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "28.904.80"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  -0.93%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.077.06"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  +0.49%  "
# Row 21
$ws.Range("E21").Value = "  -0.27%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "223.96"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  -0.84%  "
# Row 23
$ws.Range("E23").Value = "  -0.01%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.118"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  -0.03%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.001"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  +0.06%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "157.97"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  -1.65%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.467"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  -1.91%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1369"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  -1.56%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "17.86"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -0.23%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.505"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  +0.08%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.092"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  -0.41%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.023"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -0.06%  "
# Row 33
$ws.Range("E33").Value = "  +1.57%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05293"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  -0.02%  "
# Row 35
$ws.Range("E35").Value = "  -1.70%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7381"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  -2.05%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.152"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  +1.97%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.654"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  -0.90%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.291.65"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  -0.74%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.751"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  +1.20%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01781"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  -0.81%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.290"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  +5.86%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8961"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  -2.65%  "
# Row 44
$ws.Range("E44").Value = "  -0.11%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "102.75"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  +0.46%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.976.19"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  +0.28%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5136"
$ws.Range("D47").NumberFormat = "General"
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "63.92"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  +1.04%  "
# Row 49
$ws.Range("E49").Value = "  -4.85%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.728"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  -2.08%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07268"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  -15.15%  "
